$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F-SW-SD-03")
$ws.Name = "F-SW-SD-08"
